$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.667.11'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '2.518.76'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'318.88"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Range('E5').Value = '  +4.46%  '
$ws.Range('D6').Value = "'95.57"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').Value = '  +0.71%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -2.22%  '
$ws.Range('D10').Value = "'36.15"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Range('E10').Value = '  -2.95%  '
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('D12').Value = "'7.59"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Range('E12').Value = '  -2.22%  '
$ws.Range('D13').Value = "'0.113"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Range('E13').Value = '  -3.19%  '
$ws.Range('D14').Value = '2.905.20'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').Value = "'15.49"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Range('E15').Value = '  +2.90%  '
$ws.Range('D16').Value = '2.514.48'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').Value = "'0.858"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Range('E17').Value = '  -3.17%  '
$ws.Range('D18').Value = '42.691.17'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = "'12.96"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Range('E19').Value = '  -5.93%  '
$ws.Range('D20').Value = '0.0₃0971'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').Value = "'6.55"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').Value = "'71.33"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Range('E22').Value = '  -0.69%  '
$ws.Range('D23').Value = "'251.71"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('E25').Value = '  -3.22%  '
$ws.Range('D26').Value = "'26.96"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Range('E26').Value = '  -3.51%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  +12.87%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = "'38.88"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Range('E29').Value = '  +1.61%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = "'10.10"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('D32').Value = "'155.43"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('D33').Value = "'19.37"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Range('E33').Value = '  +1.94%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  -3.83%  '
$ws.Range('E36').Value = '  -3.01%  '
$ws.Range('E37').Value = '  -5.06%  '
$ws.Range('D39').Value = "'24.29"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Range('E39').Value = '  -7.64%  '
$ws.Range('D40').Value = "'0.119"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = "'3.86"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('E43').Value = '  -3.08%  '
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('D46').Value = '2.036.76'
$ws.Range('E46').Value = '  -2.94%  '
$ws.Range('D47').Value = "'84.18"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Range('E47').Value = '  -2.86%  '
$ws.Range('D48').Value = "'8.82"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Range('E48').Value = '  -3.23%  '
$ws.Range('D49').Value = '2.761.65'
$ws.Range('E49').Value = '  -1.48%  '
$ws.Range('D50').Value = "'73.27"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Range('E50').Value = '  -3.07%  '
$ws.Range('E51').Value = '  -0.93%  '
